$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.890.18"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.592.63"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.64"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.39"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.64"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.23"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.058.93"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.767.99"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.711.37"
$ws.Range("E17").Value = "  +7.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.75"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.67"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.734.83"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.60"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.00"
$ws.Range("E28").Value = "  +10.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.44"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.93"
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0824"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "458.19"
$ws.Range("E34").Value = "  +12.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.50"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.408"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.24"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.56"
$ws.Range("E40").Value = "  +5.59%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "159.00"
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.80"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.634"
$ws.Range("E45").Value = "  +5.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.70"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0546"
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0973"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0238"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.62"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.74"
$ws.Range("E51").Value = "  +0.29%  "
